$wb = $excel.ActiveWorkbook

$wsDBD = $wb.Worksheets.Item("DBD")
$wsDBS = $wb.Worksheets.Item("DBS")

# Update the "L8201維護 : xxx" labels in column G (DBD sheet) to use a
# full-width colon without surrounding spaces: "L8201維護：xxx"
$wsDBD.Range("G10").Value = "L8201維護：金額合計超過"
$wsDBD.Range("G11").Value = "L8201維護：次數 "
$wsDBD.Range("G12").Value = "L8201維護：單筆起始金額 "
$wsDBD.Range("G13").Value = "L8201維護：單筆迄止金額 "
$wsDBD.Range("G14").Value = "L8201維護：金額合計超過"
$wsDBD.Range("G15").Value = "L8201維護：統計期間天數"
$wsDBD.Range("G16").Value = "L8201維護：統計期間天數"

# Update the selection remembered on the DBS sheet before switching away
$wsDBS.Range("A24").Select() | Out-Null

# Make DBD the active sheet/tab, with its own remembered selection
$wsDBD.Activate() | Out-Null
$wsDBD.Range("G18").Select() | Out-Null
